$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("costs_inpt")

# Wrap the raw cost-conversion formulas in ROUND(...,0) so the stored
# (and displayed) values become whole numbers instead of long decimals.
$ws.Range("C2").Formula = "=ROUND(B2*((12358-6551)/3.92)/9454,0)"
$ws.Range("C3").Formula = "=ROUND(B3*((12358-6551)/3.92)/9454,0)"
$ws.Range("C4").Formula = "=ROUND(B4*((12358-6551)/3.92)/9454,0)"

# Move the active selection from C4 to C5, matching the author's last
# saved cursor position.
$ws.Range("C5").Select()
